# First commit for the assignment
# Replace the sample login test data on the "LoginData" sheet with new
# values, drop the 6th (now unused) row, and update the cell selections
# on both sheets.

$wb = $excel.ActiveWorkbook

$loginSheet = $wb.Worksheets.Item("LoginData")
$sheet1 = $wb.Worksheets.Item("Sheet1")

# New data for rows 2-5 (row 1 headers stay: firstname / lastname / postcode)
$loginSheet.Range("A2").Value = "Duleeka"
$loginSheet.Range("B2").Value = "Munasinghe"
$loginSheet.Range("C2").Value = "500128S"

$loginSheet.Range("C3").Value = "a2"
$loginSheet.Range("C4").Value = "ad3"
$loginSheet.Range("C5").Value = "bdge334"

$loginSheet.Range("A3").Value = "Kumara"
$loginSheet.Range("B3").Value = "Sangakkara"

$loginSheet.Range("A4").Value = "Aravinda"
$loginSheet.Range("B4").Value = "Silva"

$loginSheet.Range("A5").Value = "Virath"
$loginSheet.Range("B5").Value = "Kholi"

# Row 6 no longer exists in the updated data set
$loginSheet.Range("A6:C6").Clear()

# Update the remembered selection on each sheet (select the non-active
# sheet first, then re-select the LoginData sheet last so it remains the
# active tab, matching the original workbook's tab selection)
$sheet1.Range("A1:C3").Select()
$loginSheet.Range("C7").Select()
